$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B11").Value = "0b044b1b04067b410ef832fd95f1b939"
$ws.Range("B81").Value = "bc13665eac829680b6a0efac910209a9"
$ws.Range("B99").Value = "8fa707fa36d3caa430d7ac46f48d4f9c"
$ws.Range("B110").Value = "65de952e58dc1722949c67d995b7045e"
$ws.Range("B157").Value = "581e190fc1c856700a3b894d77ebe8f1"
$ws.Range("B282").Value = "dab16dfe6b5911248db88051126e75e2"
$ws.Range("B400").Value = "ca9a0ce7200f67ff0f489c634cd576bf"
$ws.Range("B404").Value = "aec11b26aac47ff6bdcac8864b6f5cbf"
$ws.Range("B409").Value = "363b8da5a43db16b69f56ba299c69d4a"
$ws.Range("B519").Value = "caf0902acd5e4ab007abd4dbb31c1a66"
$ws.Range("B547").Value = "17061e37991d1570129d34954743815d"
$ws.Range("B636").Value = "0881fa53f454181668a3a466c4556f0c"
$ws.Range("B745").Value = "2802ab1063279d54146223f696f20eb3"
$ws.Range("B748").Value = "fb6579275369feca2249f6a62946d497"
$ws.Range("B753").Value = "45cce2fdc22e2cfd7fa5302a2e549dab"
$ws.Range("B782").Value = "27b49a6dc48a01aeb632fc181b969190"
$ws.Range("B825").Value = "7e88dd68aa0fc5170b0ffca95c658e22"
$ws.Range("B829").Value = "ccde100379ce7a959dfda14c2b5d2d33"
$ws.Range("B906").Value = "2466887cd7691698e9f27cf8f11429a7"
$ws.Range("B942").Value = "0739e4252751d56b83824b70b671b54d"
$ws.Range("B965").Value = "164564ca6182282ff0c3c6b63f6c25c6"
